$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D3").Value = -7.713000000000001
$ws.Range("E3").Value = 16.885
$ws.Range("C12").Value = -10.771
$ws.Range("D14").Value = -7.597999999999999
$ws.Range("E20").Value = 16.304
$ws.Range("E25").Value = 17.014
$ws.Range("D26").Value = -7.763000000000001
$ws.Range("C27").Value = -12.771
$ws.Range("E30").Value = 16.305
$ws.Range("D31").Value = -8.123000000000001
$ws.Range("C32").Value = -13.604
$ws.Range("D35").Value = -7.935
$ws.Range("C36").Value = -12.721
$ws.Range("D37").Value = -7.712999999999999
$ws.Range("C38").Value = -12.703
$ws.Range("E44").Value = 16.613
$ws.Range("D45").Value = -7.57
$ws.Range("C46").Value = -13.736
$ws.Range("E47").Value = 16.148
$ws.Range("D52").Value = -7.281000000000001
$ws.Range("C54").Value = -12.802
$ws.Range("C55").Value = -13.391
$ws.Range("C56").Value = -13.364
$ws.Range("D57").Value = -8.184000000000001
$ws.Range("E58").Value = 16.576
$ws.Range("C67").Value = -11.661
$ws.Range("C69").Value = -10.751
$ws.Range("C72").Value = -11.555
$ws.Range("E78").Value = 16.391
$ws.Range("D81").Value = -7.085000000000001
$ws.Range("C83").Value = -13.392
$ws.Range("D83").Value = -8.334
$ws.Range("E84").Value = 16.256
$ws.Range("C86").Value = -13.846
$ws.Range("E89").Value = 17.36199999999999
$ws.Range("C91").Value = -11.308
$ws.Range("E91").Value = 17.076
$ws.Range("E92").Value = 16.883
$ws.Range("C93").Value = -11.979
$ws.Range("E96").Value = 16.283
$ws.Range("C99").Value = -12.635
$ws.Range("D100").Value = -8.005000000000001
$ws.Range("D102").Value = -7.695
$ws.Range("E102").Value = 16.425
